# Update cryptos price/volume table with latest scrape (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.301.67"
$ws.Range("E2").Value = "  -0.18%  "
$ws.Range("D3").Value = "2.981.72"
$ws.Range("E3").Value = "  +1.91%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'384.01"
$ws.Range("E5").Value = "  +1.67%  "
$ws.Range("D6").Value = "'102.47"
$ws.Range("E6").Value = "  -0.67%  "
$ws.Range("D7").Value = "'0.540"
$ws.Range("E7").Value = "  -0.49%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  +1.12%  "
$ws.Range("D10").Value = "'36.74"
$ws.Range("E10").Value = "  -1.01%  "
$ws.Range("E11").Value = "  +0.14%  "
$ws.Range("E12").Value = "  +0.88%  "
$ws.Range("D13").Value = "3.456.05"
$ws.Range("E13").Value = "  +2.26%  "
$ws.Range("D14").Value = "'18.18"
$ws.Range("E14").Value = "  -0.89%  "
$ws.Range("D15").Value = "'7.49"
$ws.Range("E15").Value = "  +1.90%  "
$ws.Range("D16").Value = "2.979.82"
$ws.Range("E16").Value = "  +2.24%  "
$ws.Range("D17").Value = "'0.999"
$ws.Range("E17").Value = "  +7.38%  "
$ws.Range("D18").Value = "51.254.39"
$ws.Range("E18").Value = "  -0.14%  "
$ws.Range("D19").Value = "'3.24"
$ws.Range("E19").Value = "  -5.14%  "
$ws.Range("D20").Value = "'7.37"
$ws.Range("E20").Value = "  -0.16%  "
$ws.Range("D21").Value = "'12.76"
$ws.Range("E21").Value = "  -1.71%  "
$ws.Range("D22").Value = "0.0₃0957"
$ws.Range("E22").Value = "  +1.25%  "
$ws.Range("D23").Value = "'68.85"
$ws.Range("E23").Value = "  +0.74%  "
$ws.Range("D24").Value = "'262.13"
$ws.Range("E24").Value = "  +0.23%  "
$ws.Range("E25").Value = "  +4.89%  "
$ws.Range("D26").Value = "'8.20"
$ws.Range("E26").Value = "  +13.63%  "
$ws.Range("D27").Value = "'7.59"
$ws.Range("E27").Value = "  +11.08%  "
$ws.Range("D28").Value = "'0.169"
$ws.Range("E28").Value = "  -1.04%  "
$ws.Range("E29").Value = "  +12.22%  "
$ws.Range("D30").Value = "'4.11"
$ws.Range("E30").Value = "  -0.09%  "
$ws.Range("E31").Value = "  +0.05%  "
$ws.Range("D32").Value = "'25.85"
$ws.Range("E32").Value = "  +0.59%  "
$ws.Range("E33").Value = "  +0.48%  "
$ws.Range("D34").Value = "'34.45"
$ws.Range("E34").Value = "  +0.96%  "
$ws.Range("D35").Value = "'50.88"
$ws.Range("E35").Value = "  -0.41%  "
$ws.Range("D36").Value = "'2.05"
$ws.Range("E36").Value = "  -2.59%  "
$ws.Range("D37").Value = "'0.0448"
$ws.Range("E37").Value = "  +6.42%  "
$ws.Range("E38").Value = "  +0.08%  "
$ws.Range("D39").Value = "'3.00"
$ws.Range("E39").Value = "  -1.07%  "
$ws.Range("D40").Value = "'17.05"
$ws.Range("E40").Value = "  +0.63%  "
$ws.Range("D41").Value = "'2.58"
$ws.Range("E41").Value = "  +0.66%  "
$ws.Range("E42").Value = "  +1.33%  "
$ws.Range("E43").Value = "  -0.85%  "
$ws.Range("D44").Value = "'122.27"
$ws.Range("E44").Value = "  -1.71%  "
$ws.Range("D45").Value = "'21.57"
$ws.Range("E45").Value = "  -0.13%  "
$ws.Range("E46").Value = "  +0.16%  "
$ws.Range("E47").Value = "  +2.01%  "
$ws.Range("E48").Value = "  +1.82%  "
$ws.Range("D49").Value = "2.033.01"
$ws.Range("E49").Value = "  +0.28%  "
$ws.Range("D50").Value = "'3.26"
$ws.Range("E50").Value = "  +3.02%  "
$ws.Range("D51").Value = "'0.0336"
$ws.Range("E51").Value = "  +5.48%  "
